$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 3.6
$ws.Range("W2").Value = 1.24
$ws.Range("J3").Value = 2.68
$ws.Range("L3").Value = 1.61
$ws.Range("Q3").Value = 3.2
$ws.Range("AB4").Value = 9.4
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AD4").Value = 20
$ws.Range("AH4").Value = 20
$ws.Range("AJ4").Value = 21
$ws.Range("AN4").Value = 12.5
$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 4.5
$ws.Range("I4").Value = 5.3
$ws.Range("J4").Value = 3.65
$ws.Range("K4").Value = 4.2
$ws.Range("P4").Value = 1.97
$ws.Range("Q4").Value = 1.84
$ws.Range("S4").Value = 3.15
$ws.Range("T4").Value = 1.79
$ws.Range("U4").Value = 2.04
$ws.Range("V4").Value = 1.23
$ws.Range("W4").Value = 2.1
$ws.Range("X4").Value = 18
$ws.Range("Y4").Value = 18.5
$ws.Range("AA5").Value = 190
$ws.Range("AM5").Value = 160
$ws.Range("F5").Value = 1.69
$ws.Range("I5").Value = 6.6
$ws.Range("N5").Value = 3.45
$ws.Range("R5").Value = 1.31
$ws.Range("U5").Value = 1.86
$ws.Range("AF6").Value = 8.6
$ws.Range("U6").Value = 1.74
$ws.Range("W6").Value = 4.7
$ws.Range("AE7").Value = 16.5
$ws.Range("F7").Value = 17
$ws.Range("H7").Value = 1.16
$ws.Range("J7").Value = 8
$ws.Range("T7").Value = 2.26
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 6
$ws.Range("Z7").Value = 9.4
$ws.Range("AA8").Value = 70
$ws.Range("AJ8").Value = 60
$ws.Range("AK8").Value = 55
$ws.Range("F8").Value = 2.9
$ws.Range("G8").Value = 3.2
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 2.72
$ws.Range("K8").Value = 2.88
$ws.Range("T8").Value = 2.32
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 1.42
$ws.Range("W8").Value = 1.46
$ws.Range("Y8").Value = 8.4
$ws.Range("AC9").Value = 9.4
$ws.Range("AF9").Value = 11
$ws.Range("F9").Value = 1.79
$ws.Range("G9").Value = 1.87
$ws.Range("H9").Value = 5.8
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 3.2
$ws.Range("K9").Value = 3.5
$ws.Range("M9").Value = 1.12
$ws.Range("N9").Value = 2.72
$ws.Range("O9").Value = 1.51
$ws.Range("P9").Value = 1.56
$ws.Range("Q9").Value = 2.48
$ws.Range("R9").Value = 1.2
$ws.Range("S9").Value = 5.1
$ws.Range("T9").Value = 2.26
$ws.Range("U9").Value = 1.69
$ws.Range("W9").Value = 2.14
$ws.Range("X9").Value = 11
$ws.Range("Y9").Value = 18.5
$ws.Range("AA10").Value = 150
$ws.Range("AD10").Value = 26
$ws.Range("F10").Value = 1.71
$ws.Range("G10").Value = 1.83
$ws.Range("H10").Value = 5.4
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 4
$ws.Range("N10").Value = 3.7
$ws.Range("P10").Value = 1.92
$ws.Range("T10").Value = 1.87
$ws.Range("U10").Value = 1.98
$ws.Range("V10").Value = 1.18
$ws.Range("W10").Value = 2.2
$ws.Range("X10").Value = 1000
$ws.Range("Z10").Value = 60
$ws.Range("J11").Value = 2.84
$ws.Range("AE12").Value = 48
$ws.Range("F12").Value = 2.26
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 3.7
$ws.Range("K12").Value = 3.4
$ws.Range("V12").Value = 1.35
$ws.Range("W12").Value = 1.77
$ws.Range("G13").Value = 2.5
$ws.Range("W13").Value = 1.7
$ws.Range("AA14").Value = 210
$ws.Range("AB14").Value = 7.4
$ws.Range("AD14").Value = 26
$ws.Range("AF14").Value = 9.199999999999999
$ws.Range("AG14").Value = 10
$ws.Range("AI14").Value = 130
$ws.Range("AJ14").Value = 16.5
$ws.Range("AK14").Value = 19
$ws.Range("AM14").Value = 210
$ws.Range("AN14").Value = 12
$ws.Range("AO14").Value = 180
$ws.Range("F14").Value = 1.66
$ws.Range("L14").Value = 1.44
$ws.Range("N14").Value = 3.35
$ws.Range("O14").Value = 1.37
$ws.Range("P14").Value = 1.84
$ws.Range("Q14").Value = 2.08
$ws.Range("R14").Value = 1.32
$ws.Range("S14").Value = 3.85
$ws.Range("T14").Value = 2.06
$ws.Range("U14").Value = 1.86
$ws.Range("W14").Value = 2.38
$ws.Range("Z14").Value = 50
$ws.Range("AA15").Value = 470
$ws.Range("AC15").Value = 970
$ws.Range("AJ15").Value = 1000
$ws.Range("AN15").Value = 9.199999999999999
$ws.Range("AO15").Value = 390
$ws.Range("F15").Value = 1.46
$ws.Range("J15").Value = 4.4
$ws.Range("K15").Value = 4.7
$ws.Range("L15").Value = 1.42
$ws.Range("P15").Value = 1.83
$ws.Range("Q15").Value = 2.08
$ws.Range("T15").Value = 2.28
$ws.Range("U15").Value = 1.69
$ws.Range("AA16").Value = 580
$ws.Range("AC16").Value = 10.5
$ws.Range("AD16").Value = 970
$ws.Range("AI16").Value = 260
$ws.Range("F16").Value = 1.5
$ws.Range("G16").Value = 1.55
$ws.Range("W16").Value = 2.8
